# Fruta / hortaliza, semanal
# Insert a new weekly price record for "Vega Monumental Concepción - Mandarina"
# above the existing row 212, shifting the remaining rows (212-220) down to
# (213-221) and growing the used range to A1:T221.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 212 (pushes old rows 212..220 -> 213..221)
$ws.Rows.Item(212).Insert()

# Populate the newly inserted row with the new weekly observation
$ws.Range("A212").Value = 11
$ws.Range("B212").Value = "Vega Monumental Concepción"
$ws.Range("C212").Value = "Bíobío"
$ws.Range("D212").Value = 45106
$ws.Range("E212").Value = 8
$ws.Range("F212").Value = "Fruta"
$ws.Range("G212").Value = 100102
$ws.Range("H212").Value = "Cítricos"
$ws.Range("I212").Value = 100102004
$ws.Range("J212").Value = "Mandarina"
$ws.Range("K212").Value = "Murcott"
$ws.Range("L212").Value = "Primera"
$ws.Range("M212").Value = 150
$ws.Range("N212").Value = 9000
$ws.Range("O212").Value = 10000
$ws.Range("P212").Value = 9467
$ws.Range("Q212").Value = "$/bandeja 10 kilos"
$ws.Range("R212").Value = "Región de O'Higgins"
$ws.Range("S212").Value = 947
$ws.Range("T212").Value = 10
